# Curve Prep - Ramer Douglas Peucker
# Remove the data rows that the RDP simplification drops, shifting the
# remaining rows up (Excel's normal EntireRow.Delete behaviour).
# Rows are deleted from bottom to top so row numbers of not-yet-deleted
# rows stay stable while iterating.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToDelete = @(50, 48, 46, 42, 40, 38, 36, 34, 30, 27, 25, 23, 21, 5, 3)

foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}

Write-Host "Deleted $($rowsToDelete.Count) rows; new dimension now A1:B39"
